$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions refresh).
$ws.Range('D2').Value = '63.844.65'
$ws.Range('E2').Value = '  -1.30%  '
$ws.Range('D3').Value = '3.135.67'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '598.79'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '139.57'
$ws.Range('E6').Value = '  -4.00%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = '3.124.59'
$ws.Range('E8').Value = '  -1.09%  '
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  -2.64%  '
$ws.Range('D11').Value = '5.31'
$ws.Range('E11').Value = '  -2.55%  '
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('E13').Value = '  -2.33%  '
$ws.Range('D14').Value = '34.49'
$ws.Range('E14').Value = '  -2.96%  '
$ws.Range('D15').Value = '3.650.76'
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '63.859.72'
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').Value = '3.136.17'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').Value = '6.76'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').Value = '481.67'
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = '14.51'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('E22').Value = '  -2.06%  '
$ws.Range('D23').Value = '7.67'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '87.54'
$ws.Range('E24').Value = '  +4.26%  '
$ws.Range('E25').Value = '  -5.06%  '
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('E28').Value = '  -6.18%  '
$ws.Range('D29').Value = '6.94'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').Value = '2.05'
$ws.Range('E30').Value = '  -2.49%  '
$ws.Range('D31').Value = '27.18'
$ws.Range('E31').Value = '  +2.55%  '
$ws.Range('E32').Value = '  -7.27%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '2.58'
$ws.Range('E34').Value = '  -3.63%  '
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('E36').Value = '  +0.09%  '
$ws.Range('D37').Value = '52.52'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('D38').Value = '0.0₃0732'
$ws.Range('E38').Value = '  -6.84%  '
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').Value = '429.28'
$ws.Range('E40').Value = '  -6.73%  '
$ws.Range('E41').Value = '  -11.15%  '
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').Value = '8.29'
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = '2.877.59'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('E45').Value = '  -3.30%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = '2.37'
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '2.15'
$ws.Range('E47').Value = '  -7.23%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').Value = '25.51'
$ws.Range('E50').Value = '  -3.83%  '
$ws.Range('D51').Value = '120.44'
$ws.Range('E51').Value = '  +0.51%  '
